$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 240, shifting rows 240:276 down to 241:277.
$ws.Rows.Item(240).Insert()

# Populate the newly inserted row 240 with the new record's data.
$ws.Cells.Item(240, 1).Value = 7
$ws.Cells.Item(240, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(240, 3).Value = "Ñuble"
$ws.Cells.Item(240, 4).Value = 44951
$ws.Cells.Item(240, 5).Value = 16
$ws.Cells.Item(240, 6).Value = "Fruta"
$ws.Cells.Item(240, 7).Value = 100108
$ws.Cells.Item(240, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(240, 9).Value = 100108005
$ws.Cells.Item(240, 10).Value = "Piña"
$ws.Cells.Item(240, 11).Value = "Caramelo"
$ws.Cells.Item(240, 12).Value = "Segunda"
$ws.Cells.Item(240, 13).Value = 60
$ws.Cells.Item(240, 14).Value = 1600
$ws.Cells.Item(240, 15).Value = 17000
$ws.Cells.Item(240, 16).Value = 9300
$ws.Cells.Item(240, 17).Value = "`$/caja 14 unidades"
$ws.Cells.Item(240, 18).Value = "Ecuador"
$ws.Cells.Item(240, 19).Value = 664
$ws.Cells.Item(240, 20).Value = 14
